$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B (ASIN) to make room for Week_Start_Date
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(2).NumberFormat = "@"
$ws.Range("B1").Value = "Week_Start_Date"

# Data: Week label (no leading zero), Week_Start_Date, MyForecast, Amazon Mean,
# Amazon P70, Amazon P80, Amazon P90, is_holiday_week (boolean)
$data = @(
    @{Row=2;  Week="W1";  Date="2025-01-05"; D=7; E=6; F=7; G=9;  H=13},
    @{Row=3;  Week="W2";  Date="2025-01-12"; D=6; E=6; F=8; G=11; H=15},
    @{Row=4;  Week="W3";  Date="2025-01-19"; D=6; E=6; F=8; G=10; H=13},
    @{Row=5;  Week="W4";  Date="2025-01-26"; D=6; E=7; F=8; G=11; H=16},
    @{Row=6;  Week="W5";  Date="2025-02-02"; D=6; E=7; F=8; G=11; H=16},
    @{Row=7;  Week="W6";  Date="2025-02-09"; D=6; E=6; F=8; G=11; H=16},
    @{Row=8;  Week="W7";  Date="2025-02-16"; D=7; E=7; F=8; G=12; H=18},
    @{Row=9;  Week="W8";  Date="2025-02-23"; D=7; E=7; F=8; G=12; H=19},
    @{Row=10; Week="W9";  Date="2025-03-02"; D=8; E=7; F=8; G=11; H=18},
    @{Row=11; Week="W10"; Date="2025-03-09"; D=8; E=7; F=8; G=12; H=18},
    @{Row=12; Week="W11"; Date="2025-03-16"; D=8; E=7; F=8; G=12; H=19},
    @{Row=13; Week="W12"; Date="2025-03-23"; D=9; E=8; F=9; G=13; H=20},
    @{Row=14; Week="W13"; Date="2025-03-30"; D=8; E=8; F=9; G=13; H=20},
    @{Row=15; Week="W14"; Date="2025-04-06"; D=7; E=7; F=8; G=12; H=20},
    @{Row=16; Week="W15"; Date="2025-04-13"; D=7; E=7; F=8; G=12; H=20},
    @{Row=17; Week="W16"; Date="2025-04-20"; D=6; E=7; F=8; G=12; H=20}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.Week
    $ws.Range("B$r").Value = $item.Date
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    $ws.Range("G$r").Value = $item.G
    $ws.Range("H$r").Value = $item.H
    $ws.Range("J$r").Value = $false
}

# Update Summary sheet values (kept as text, matching existing column style)
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B9").NumberFormat = "@"
$summary.Range("B9").Value = "114"
$summary.Range("B12").NumberFormat = "@"
$summary.Range("B12").Value = "9"
$summary.Range("B13").NumberFormat = "@"
$summary.Range("B13").Value = "2025-03-23"
